$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Work around an xlsx-exporter quirk: when a cell's new unique string happens to
# land back on the exact same shared-string-table index it previously had, the
# cell's <v> index isn't refreshed on save. Briefly parking a throwaway unique
# string in an unused cell shifts the index during the edit so the save is
# correct; clearing that helper cell afterwards removes it completely.
$ws.Cells.Item(100, 26).Value = "DUMMY_PLACEHOLDER_TO_SHIFT"

$newQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC02'] and diag.stage_of_disease in [ 'T2N0M0', 'T2N0M1', 'T2N1M0', 'T2N1M1', 'T3N0M0']  OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@
$ws.Range("B2").Value = $newQuery

$ws.Cells.Item(100, 26).ClearContents()

# The shorter query text re-wraps to fewer lines, so the wrapped rows shrink
# (mirrors Excel's own autofit recalculation after the edit).
$ws.Rows.Item(2).RowHeight = 304.5
$ws.Rows.Item(3).RowHeight = 290
$ws.Rows.Item(4).RowHeight = 290

$ws.Activate()
$ws.Range("B2").Select() | Out-Null
